$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Remove the first comment thread ("Is this correct, Ian?") together with
#    its in-text comment range markers (commentRangeStart/End + reference).
#    Comment.Delete() removes the markers from the body AND the <w:comment>
#    entry from comments.xml (and renumbers nothing else).
# ---------------------------------------------------------------------------
$targetComment = $null
for ($i = 1; $i -le $d.Comments.Count; $i++) {
    $c = $d.Comments.Item($i)
    if ($c.Range.Text -eq "Is this correct, Ian?") {
        $targetComment = $c
    }
}
if ($targetComment -ne $null) {
    $targetComment.Delete()
}

# ---------------------------------------------------------------------------
# 2. Word always keeps a single "_GoBack" bookmark marking the most recent
#    edit location. Since the edit above happened right after "...The D2 ",
#    move/re-create the bookmark there (adding a bookmark named "_GoBack"
#    automatically removes any pre-existing one elsewhere in the document).
# ---------------------------------------------------------------------------
$findRange = $d.Content
$found = $findRange.Find.Execute("The D2 ", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $editPoint = $findRange.End
    $d.Bookmarks.Add("_GoBack", $d.Range($editPoint, $editPoint))
}
